$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 - Potentiometre
$ws.Range("B11").Value = "Potentiometre"
$ws.Range("C11").Value = "Potentiometre 10K"
$ws.Range("D11").Value = 2
$ws.Range("E11").Value = 0.68
$ws.Range("F11").Formula = "=E11*D11+0.99"
$ws.Range("G11").Value = 10
$ws.Range("H11").Value = "https://fr.aliexpress.com/item/32996236826.html?spm=a2g0o.detail.1000060.3.2dfd75067cm24G&gps-id=pcDetailBottomMoreThisSeller&scm=1007.13339.291025.0&scm_id=1007.13339.291025.0&scm-url=1007.13339.291025.0&pvid=aa5bbdd6-37f4-4697-8922-a5d47906774f&_t=gps-id%3ApcDetailBottomMoreThisSeller%2Cscm-url%3A1007.13339.291025.0%2Cpvid%3Aaa5bbdd6-37f4-4697-8922-a5d47906774f%2Ctpp_buckets%3A668%232846%238116%232002&pdp_ext_f=%7B%22sku_id%22%3A%2266994016160%22%2C%22sceneId%22%3A%223339%22%7D&pdp_npi=2%40dis%21EUR%210.77%210.68%21%21%21%21%21%402103222316652351797365094e1aae%2166994016160%21rec&gatewayAdapt=glo2fra"
$ws.Range("I11").Value = "Color: 10 K Ohm"

# Row 12 - Accelerometre
$ws.Range("B12").Value = "Accelerometre"
$ws.Range("C12").Value = "MPU6050"
$ws.Range("D12").Value = 1
$ws.Range("E12").Formula = "=1.42+1.17"
$ws.Range("F12").Formula = "=E12*D12"
$ws.Range("G12").Value = 6
$ws.Range("H12").Value = "https://fr.aliexpress.com/item/32340949017.html?spm=a2g0o.productlist.0.0.587312adnNgBZj&algo_pvid=19376db7-79d4-41f3-91f8-c1871eddb2a7&algo_exp_id=19376db7-79d4-41f3-91f8-c1871eddb2a7-0&pdp_ext_f=%7B%22sku_id%22%3A%2210000000609322940%22%7D&pdp_npi=2%40dis%21EUR%211.62%211.42%21%21%211.71%21%21%402101e9d416652354615517308e70c3%2110000000609322940%21sea&curPageLogUid=c2sYDhCk9WsN"
$ws.Range("I12").Value = " "

# Selection
$ws.Range("B13").Select()
